# SonarLint-proposed changes: fill in the "Tool-based Code Analysis" sheet
# with the actual findings/fixes (the sheet previously only had the header
# rows and 21 empty data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tool-basedCodeAnalysis")

# Reviewer/tool-used answer (row 5, column D) on the summary block.
$ws.Range("D5").Value = '3testeri'

# --- Data rows (Crt.No already 1..21 via formulas; fill File/Issue/Before/After) ---

# Row 10
$ws.Range("C10").Value = 'Task, line 125'
$ws.Range("D10").Value = 'Refactor nextTimeAfter method to reduce its Cognitive Complexity'
$ws.Range("E10").Value = 'Codul este ingramadit si greu de inteles'
$ws.Range("F10").Value = 'Am creat metode noi care au fiecare propriul rol.'

# Row 11
$ws.Range("C11").Value = 'Task, line 233'
$ws.Range("D11").Value = 'Remove this "clone"  impl'
$ws.Range("F11").Value = 'Am adaugat un constructor care are ca parametru un obiect de tip task si salveaza atributele acestuia'

# Row 12
$ws.Range("C12").Value = 'Task, line 125'
$ws.Range("D12").Value = 'Use the opposite operator "!="'
$ws.Range("E12").Value = '!(this.interval == 0)'
$ws.Range("F12").Value = 'this.interval != 0'

# Row 13
$ws.Range("C13").Value = 'Task, line 146, 150'
$ws.Range("D13").Value = 'Remove unused method parameter "current"'
$ws.Range("E13").Value = 'current este parametru la cele 2 functii'
$ws.Range("F13").Value = 'l-am eliminat ca parametru'

# Row 14
$ws.Range("C14").Value = 'Task, line 13'
$ws.Range("D14").Value = 'Make "sdf" an instance variable'
$ws.Range("E14").Value = 'private static final SimpleDateFormat sdf'
$ws.Range("F14").Value = 'eliminarea keyword ului static duce la erori in teste'

# Row 15
$ws.Range("C15").Value = 'TaskList, line 36'
$ws.Range("D15").Value = 'Replace sout by a logger'
$ws.Range("E15").Value = 'System.out.println(getTask(i).getTitle());'
$ws.Range("F15").Value = 'logger.info(getTask(i).getTitle());'

# Row 16 (E16 keeps the leading newline/indentation of the removed code line)
$ws.Range("C16").Value = 'TaskList, line 21'
$ws.Range("D16").Value = '"Iterator" is defined int the "Iterable" interafce and can be removed from this class'
$ws.Range("E16").Value = "`n   public abstract Iterator<Task> iterator();"
$ws.Range("F16").Value = 'am eliminat aceasta linie de cod'

# --- Cosmetic touch-ups that accompanied the content edit ---

# Columns D and F had to grow to host the new, longer text.
$ws.Columns.Item(4).ColumnWidth = 23.83
$ws.Columns.Item(6).ColumnWidth = 23.92

# Rows with wrapped multi-line content grew taller.
$ws.Rows.Item(10).RowHeight = 72.5
$ws.Rows.Item(11).RowHeight = 77
$ws.Rows.Item(13).RowHeight = 29
$ws.Rows.Item(14).RowHeight = 29
$ws.Rows.Item(15).RowHeight = 29
$ws.Rows.Item(16).RowHeight = 43.5

# The reviewer scrolled down a bit and left the cursor on I15.
$ws.Activate()
$ws.Range("I15").Select()
